$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.057028082059764
$ws.Range("D2").Value = 1.054582980607437
$ws.Range("E2").Value = 1.06229614213626
$ws.Range("F2").Value = 1.071272194292041
$ws.Range("I2").Value = 1.041541081105195
$ws.Range("J2").Value = 1.062026641344852
$ws.Range("K2").Value = 1.057325791949645
$ws.Range("L2").Value = 1.065017874336167
$ws.Range("M2").Value = 1.073969809336351

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058887323096292
$ws.Range("D3").Value = 1.056009820986789
$ws.Range("E3").Value = 1.064012694455673
$ws.Range("F3").Value = 1.073179524846739
$ws.Range("I3").Value = 1.042037246872855
$ws.Range("J3").Value = 1.063533812675358
$ws.Range("K3").Value = 1.058564337102249
$ws.Range("L3").Value = 1.06654694832461
$ws.Range("M3").Value = 1.075690969176053

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060086506140556
$ws.Range("D4").Value = 1.056929563082407
$ws.Range("E4").Value = 1.065119628164644
$ws.Range("F4").Value = 1.074410140828672
$ws.Range("I4").Value = 1.042355048962529
$ws.Range("J4").Value = 1.064504932987288
$ws.Range("K4").Value = 1.059361726369345
$ws.Range("L4").Value = 1.067532116180817
$ws.Range("M4").Value = 1.076800699586076

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060589736302966
$ws.Range("D5").Value = 1.057315394545962
$ws.Range("E5").Value = 1.06558409414173
$ws.Range("F5").Value = 1.074926661026741
$ws.Range("I5").Value = 1.042487881101728
$ws.Range("J5").Value = 1.064912222154091
$ws.Range("K5").Value = 1.059695997014648
$ws.Range("L5").Value = 1.067945281716252
$ws.Range("M5").Value = 1.077266297409579

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060674178291912
$ws.Range("D6").Value = 1.057380129199357
$ws.Range("E6").Value = 1.065662028427667
$ws.Range("F6").Value = 1.075013338906868
$ws.Range("I6").Value = 1.042510139121102
$ws.Range("J6").Value = 1.064980551405992
$ws.Range("K6").Value = 1.059752067102848
$ws.Range("L6").Value = 1.068014595942983
$ws.Range("M6").Value = 1.077344419173302

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060093233860752
$ws.Range("D7").Value = 1.056934721815672
$ws.Range("E7").Value = 1.065125837849639
$ws.Range("F7").Value = 1.074417045833765
$ws.Range("I7").Value = 1.042356826896374
$ws.Range("J7").Value = 1.064510378993202
$ws.Range("K7").Value = 1.059366196631578
$ws.Range("L7").Value = 1.067537640820942
$ws.Range("M7").Value = 1.076806924566868

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.057657234821419
$ws.Range("D8").Value = 1.055065925320916
$ws.Range("E8").Value = 1.062877055041986
$ws.Range("F8").Value = 1.071917533755744
$ws.Range("I8").Value = 1.041709439456033
$ws.Range("J8").Value = 1.062536860137193
$ws.Range("K8").Value = 1.05774520717945
$ws.Range("L8").Value = 1.065535521520611
$ws.Range("M8").Value = 1.07455231780933

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.053334141898692
$ws.Range("D9").Value = 1.051745266633999
$ws.Range("E9").Value = 1.058884547835644
$ws.Range("F9").Value = 1.067484940558738
$ws.Range("I9").Value = 1.040543482043444
$ws.Range("J9").Value = 1.059026961010733
$ws.Range("K9").Value = 1.054857322554997
$ws.Range("L9").Value = 1.061974236987057
$ws.Range("M9").Value = 1.070548118992964

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.05043025894321
$ws.Range("D10").Value = 1.049512023025341
$ws.Range("E10").Value = 1.056201623539915
$ws.Range("F10").Value = 1.064509682149195
$ws.Range("I10").Value = 1.039748847424511
$ws.Range("J10").Value = 1.056664251418192
$ws.Range("K10").Value = 1.052910029353237
$ws.Range("L10").Value = 1.059576563544745
$ws.Range("M10").Value = 1.067856418751999

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049167382954214
$ws.Range("D11").Value = 1.048540180199646
$ws.Range("E11").Value = 1.055034583741953
$ws.Range("F11").Value = 1.063216285067726
$ws.Range("I11").Value = 1.039400561045359
$ws.Range("J11").Value = 1.055635532853797
$ws.Range("K11").Value = 1.052061409389598
$ws.Range("L11").Value = 1.058532528170124
$ws.Range("M11").Value = 1.066685341407287

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.048697447077058
$ws.Range("D12").Value = 1.048178450458764
$ws.Range("E12").Value = 1.054600271743847
$ws.Range("F12").Value = 1.062735070133362
$ws.Range("I12").Value = 1.039270552420888
$ws.Range("J12").Value = 1.055252551610407
$ws.Range("K12").Value = 1.051745361641216
$ws.Range("L12").Value = 1.058143830152938
$ws.Range("M12").Value = 1.066249493596746

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.048798288672281
$ws.Range("D13").Value = 1.04825607667494
$ws.Range("E13").Value = 1.05469347068371
$ws.Range("F13").Value = 1.062838328507408
$ws.Range("I13").Value = 1.0392984687815
$ws.Range("J13").Value = 1.055334742060854
$ws.Range("K13").Value = 1.05181319290179
$ws.Range("L13").Value = 1.058227248124881
$ws.Range("M13").Value = 1.066343023613214

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049128555341976
$ws.Range("D14").Value = 1.048510294788874
$ws.Range("E14").Value = 1.054998700267168
$ws.Range("F14").Value = 1.06317652394525
$ws.Range("I14").Value = 1.03938982758478
$ws.Range("J14").Value = 1.055603893359486
$ws.Range("K14").Value = 1.052035301865343
$ws.Range("L14").Value = 1.05850041668969
$ws.Range("M14").Value = 1.066649331720316

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049331930484707
$ws.Range("D15").Value = 1.048666827827674
$ws.Range("E15").Value = 1.055186652749112
$ws.Range("F15").Value = 1.063384791927598
$ws.Range("I15").Value = 1.039446031772953
$ws.Range("J15").Value = 1.055769610538447
$ws.Range("K15").Value = 1.05217203970484
$ws.Range("L15").Value = 1.058668605380345
$ws.Range("M15").Value = 1.066837943922045

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050513954385573
$ws.Range("D16").Value = 1.049576417676735
$ws.Range("E16").Value = 1.056278962188932
$ws.Range("F16").Value = 1.064595411286393
$ws.Range("I16").Value = 1.039771872792072
$ws.Range("J16").Value = 1.056732403299951
$ws.Range("K16").Value = 1.052966233647403
$ws.Range("L16").Value = 1.059645728150368
$ws.Range("M16").Value = 1.067934020321291

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05125392334935
$ws.Range("D17").Value = 1.050145672743244
$ws.Range("E17").Value = 1.056962699829271
$ws.Range("F17").Value = 1.065353420486849
$ws.Range("I17").Value = 1.039975133091841
$ws.Range("J17").Value = 1.057334810748348
$ws.Range("K17").Value = 1.053462945173146
$ws.Range("L17").Value = 1.060257076419969
$ws.Range("M17").Value = 1.068620057004896

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051685008116539
$ws.Range("D18").Value = 1.050477244438251
$ws.Range("E18").Value = 1.057361001077912
$ws.Range("F18").Value = 1.065795064960649
$ws.Range("I18").Value = 1.040093286267253
$ws.Range("J18").Value = 1.05768564110799
$ws.Range("K18").Value = 1.053752145917169
$ws.Range("L18").Value = 1.060613104763187
$ws.Range("M18").Value = 1.069019676358535

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051831908307284
$ws.Range("D19").Value = 1.050590223455657
$ws.Range("E19").Value = 1.057496725462484
$ws.Range("F19").Value = 1.065945572038393
$ws.Range("I19").Value = 1.040133504968791
$ws.Range("J19").Value = 1.057805173588994
$ws.Range("K19").Value = 1.053850667746662
$ws.Range("L19").Value = 1.060734406755159
$ws.Range("M19").Value = 1.069155846312812

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051174586331476
$ws.Range("D20").Value = 1.050084645303677
$ws.Range("E20").Value = 1.056889394306839
$ws.Range("F20").Value = 1.065272144080573
$ws.Range("I20").Value = 1.039953367145537
$ws.Range("J20").Value = 1.057270234517172
$ws.Range("K20").Value = 1.05340970691237
$ws.Range("L20").Value = 1.060191542692243
$ws.Range("M20").Value = 1.06854650718721

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049031323636169
$ws.Range("D21").Value = 1.048435454574414
$ws.Range("E21").Value = 1.054908840623337
$ws.Range("F21").Value = 1.06307695581827
$ws.Range("I21").Value = 1.039362942413178
$ws.Range("J21").Value = 1.055524659175706
$ws.Range("K21").Value = 1.051969919429258
$ws.Range("L21").Value = 1.058420000284767
$ws.Range("M21").Value = 1.066559155477641

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047678853340303
$ws.Range("D22").Value = 1.047394229069326
$ws.Range("E22").Value = 1.053658823709599
$ws.Range("F22").Value = 1.061692172479676
$ws.Range("I22").Value = 1.038988014677525
$ws.Range("J22").Value = 1.05442210640793
$ws.Range("K22").Value = 1.05105984273556
$ws.Range("L22").Value = 1.057300961899153
$ws.Range("M22").Value = 1.065304655376599

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.048396297804224
$ws.Range("D23").Value = 1.047946617337867
$ws.Range("E23").Value = 1.054321940612509
$ws.Range("F23").Value = 1.062426714979055
$ws.Range("I23").Value = 1.039187124754817
$ws.Range("J23").Value = 1.055007075188573
$ws.Range("K23").Value = 1.051542754664304
$ws.Range("L23").Value = 1.057894685366541
$ws.Range("M23").Value = 1.065970168849646

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051210436941368
$ws.Range("D24").Value = 1.050112222413635
$ws.Range("E24").Value = 1.056922519496479
$ws.Range("F24").Value = 1.065308870901851
$ws.Range("I24").Value = 1.039963203491385
$ws.Range("J24").Value = 1.057299415415909
$ws.Range("K24").Value = 1.053433764605471
$ws.Range("L24").Value = 1.060221156295085
$ws.Range("M24").Value = 1.068579742834062

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.054455516298977
$ws.Range("D25").Value = 1.052607100643385
$ws.Range("E25").Value = 1.059920365656841
$ws.Range("F25").Value = 1.06863433698814
$ws.Range("I25").Value = 1.040847934668216
$ws.Range("J25").Value = 1.059938290402981
$ws.Range("K25").Value = 1.055607730173483
$ws.Range("L25").Value = 1.062898972191632
$ws.Range("M25").Value = 1.071587130839684
